$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for the new "THREE WHEELER(T)" data row: push the
#    existing "TWO WHEELER(NT)" row down from row 5 to row 6.
# ------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ------------------------------------------------------------------
# 2. Break apart the merges that need to grow before we touch the
#    cells underneath them.
# ------------------------------------------------------------------
$ws.Range("A1:K1").UnMerge()
$ws.Range("C2:J2").UnMerge()
$ws.Range("K2:K3").UnMerge()

# ------------------------------------------------------------------
# 3. Row 2 / Row 4 headers: TOTAL moves from column K to column N,
#    and K/L/M become the new SEP / OCT / NOV month columns.
# ------------------------------------------------------------------
$ws.Cells.Item(2, 14).Value = $ws.Cells.Item(2, 11).Text
$ws.Cells.Item(2, 11).Value = "'"

$ws.Cells.Item(4, 14).Value = $ws.Cells.Item(4, 11).Text
$ws.Cells.Item(4, 11).Value = "SEP"
$ws.Cells.Item(4, 12).Value = "OCT"
$ws.Cells.Item(4, 13).Value = "NOV"

# ------------------------------------------------------------------
# 4. Re-create the merges at their new extents.
# ------------------------------------------------------------------
$ws.Range("A1:N1").Merge()
$ws.Range("C2:M2").Merge()
$ws.Range("N2:N3").Merge()

# ------------------------------------------------------------------
# 5. Row 5 data: "1", "THREE WHEELER(T)", all-zero Jan..Sep, 1 in
#    Oct, 0 in Nov, Total 1.
# ------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = "'1"
$ws.Cells.Item(5, 2).Value = "THREE WHEELER(T)"
$ws.Cells.Item(5, 3).Value = "'0"
$ws.Cells.Item(5, 4).Value = "'0"
$ws.Cells.Item(5, 5).Value = "'0"
$ws.Cells.Item(5, 6).Value = "'0"
$ws.Cells.Item(5, 7).Value = "'0"
$ws.Cells.Item(5, 8).Value = "'0"
$ws.Cells.Item(5, 9).Value = "'0"
$ws.Cells.Item(5, 10).Value = "'0"
$ws.Cells.Item(5, 11).Value = "'0"
$ws.Cells.Item(5, 12).Value = "'1"
$ws.Cells.Item(5, 13).Value = "'0"
$ws.Cells.Item(5, 14).Value = "'1"

# ------------------------------------------------------------------
# 6. Row 6 data (was row 5): "2", "TWO WHEELER(NT)", extended with
#    Sep/Oct/Nov and a refreshed Total of 7.
# ------------------------------------------------------------------
$ws.Cells.Item(6, 1).Value = "'2"
$ws.Cells.Item(6, 2).Value = "TWO WHEELER(NT)"
$ws.Cells.Item(6, 11).Value = "'1"
$ws.Cells.Item(6, 12).Value = "'0"
$ws.Cells.Item(6, 13).Value = "'0"
$ws.Cells.Item(6, 14).Value = "'7"

# ------------------------------------------------------------------
# 7. Column widths for the newly introduced / resized columns.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.39
$ws.Columns.Item(11).ColumnWidth = 3.64
$ws.Columns.Item(12).ColumnWidth = 4.05
$ws.Columns.Item(13).ColumnWidth = 4.28
$ws.Columns.Item(14).ColumnWidth = 1.38

Write-Output "edit complete"
